$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Rewrite the paragraph that begins "Soar-RL's default exploration policy
#    is ..." into the new, longer explanation about softmax / epsilon-greedy.
# ---------------------------------------------------------------------------
$targetPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -like "Soar-RL*default exploration policy is*") {
        $targetPara = $cand
        break
    }
}

$r = $targetPara.Range
# Range covering just the paragraph's text (exclude the trailing paragraph mark)
$textOnly = $d.Range($r.Start, $r.End - 1)
$textOnly.Text = ""

$ins = $d.Range($r.Start, $r.Start)

$ins.InsertAfter("When Soar is first started, the default exploration policy is ")
$ins.Collapse(0)

$ins.InsertAfter("softmax")
$ins.Font.Italic = $true
$ins.Collapse(0)

$ins.InsertAfter(".  However, the first time Soar-RL is enabled, the architecture automatically changes the exploration policy to ")
$ins.Collapse(0)

$ins.InsertAfter("epsilon-greedy")
$ins.Font.Italic = $true
$ins.Collapse(0)

$ins.InsertAfter(", ")
$ins.Collapse(0)

$ins.InsertAfter("a policy more suitable for RL agents")
$ins.Collapse(0)

$ins.InsertAfter(".  ")
$ins.Collapse(0)

$ins.InsertAfter("The default value of ")
$ins.Collapse(0)

$ins.InsertAfter("epsilon")
$ins.Font.Italic = $true
$ins.Collapse(0)

$ins.InsertAfter(" is 0.1, dictating that ")
$ins.Collapse(0)

$ins.InsertAfter("90% of the time the operator with greatest numerical preference value is chosen, while the remaining 10% of the time a random selection is made from all acceptable proposed operators.  You can change the ")
$ins.Collapse(0)

$ins.InsertAfter("epsilon")
$ins.Font.Italic = $true
$ins.Collapse(0)

$ins.InsertAfter(" value by issuing the following command:")
$ins.Collapse(0)

# ---------------------------------------------------------------------------
# 2) Italicize "epsilon" inside the "Acceptable values for epsilon are
#    numbers between 0 and 1" paragraph.
# ---------------------------------------------------------------------------
$acceptablePara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $cand = $d.Paragraphs.Item($i)
    if ($cand.Range.Text -like "Acceptable values for epsilon are numbers between 0 and 1*") {
        $acceptablePara = $cand
        break
    }
}

$ar = $acceptablePara.Range
$afind = $d.Range($ar.Start, $ar.End)
$found = $afind.Find.Execute("Acceptable values for epsilon are numbers between 0 and 1", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$startPos = $afind.Start
$afind.Text = "Acceptable values for epsilon are numbers between 0 and 1"

$italicRng = $d.Range($startPos, $startPos + 59)
$italicFind = $italicRng.Find.Execute("epsilon", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$italicRng.Font.Italic = $true

# ---------------------------------------------------------------------------
# 3) Append a new, empty paragraph at the very end of the document (just
#    before the final section properties).
# ---------------------------------------------------------------------------
$endPos = $d.Content.End
$endRng = $d.Range($endPos, $endPos)
$endRng.InsertBefore("`r")

Write-Host "Done."
